$wb = $excel.ActiveWorkbook

# Sheet 1: 大智投资组合 - append daily snapshot rows 54-76
$ws = $wb.Worksheets("大智投资组合")

$ws.Range("B54:B76").NumberFormat = "@"
$ws.Range("E54:E76").NumberFormat = "@"

$ws.Cells.Item(54, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(54, 2).Value = "600900"
$ws.Cells.Item(54, 3).Value = "长江电力"
$ws.Cells.Item(54, 4).Value = 17
$ws.Cells.Item(54, 5).Value = "202507021030"

$ws.Cells.Item(55, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(55, 2).Value = "000333"
$ws.Cells.Item(55, 3).Value = "美的集团"
$ws.Cells.Item(55, 4).Value = 3.06
$ws.Cells.Item(55, 5).Value = "202507021030"

$ws.Cells.Item(56, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(56, 2).Value = "518880"
$ws.Cells.Item(56, 3).Value = "黄金ETF"
$ws.Cells.Item(56, 4).Value = 4.88
$ws.Cells.Item(56, 5).Value = "202507021030"

$ws.Cells.Item(57, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(57, 2).Value = "510300"
$ws.Cells.Item(57, 3).Value = "沪深300ETF"
$ws.Cells.Item(57, 4).Value = 5.01
$ws.Cells.Item(57, 5).Value = "202507021030"

$ws.Cells.Item(58, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(58, 2).Value = "HK02899"
$ws.Cells.Item(58, 3).Value = "紫金矿业"
$ws.Cells.Item(58, 4).Value = 9.89
$ws.Cells.Item(58, 5).Value = "202507021030"

$ws.Cells.Item(59, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(59, 2).Value = "HK06881"
$ws.Cells.Item(59, 3).Value = "中国银河"
$ws.Cells.Item(59, 4).Value = 5.22
$ws.Cells.Item(59, 5).Value = "202507021030"

$ws.Cells.Item(60, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(60, 2).Value = "510050"
$ws.Cells.Item(60, 3).Value = "上证50ETF"
$ws.Cells.Item(60, 4).Value = 5.14
$ws.Cells.Item(60, 5).Value = "202507021030"

$ws.Cells.Item(61, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(61, 2).Value = "600085"
$ws.Cells.Item(61, 3).Value = "同仁堂"
$ws.Cells.Item(61, 4).Value = 1.98
$ws.Cells.Item(61, 5).Value = "202507021030"

$ws.Cells.Item(62, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(62, 2).Value = "600989"
$ws.Cells.Item(62, 3).Value = "宝丰能源"
$ws.Cells.Item(62, 4).Value = 1
$ws.Cells.Item(62, 5).Value = "202507021030"

$ws.Cells.Item(63, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(63, 2).Value = "601899"
$ws.Cells.Item(63, 3).Value = "紫金矿业"
$ws.Cells.Item(63, 4).Value = 9.92
$ws.Cells.Item(63, 5).Value = "202507021030"

$ws.Cells.Item(64, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(64, 2).Value = "601688"
$ws.Cells.Item(64, 3).Value = "华泰证券"
$ws.Cells.Item(64, 4).Value = 5
$ws.Cells.Item(64, 5).Value = "202507021030"

$ws.Cells.Item(65, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(65, 2).Value = "600900"
$ws.Cells.Item(65, 3).Value = "长江电力"
$ws.Cells.Item(65, 4).Value = 17
$ws.Cells.Item(65, 5).Value = "202507021326"

$ws.Cells.Item(66, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(66, 2).Value = "000333"
$ws.Cells.Item(66, 3).Value = "美的集团"
$ws.Cells.Item(66, 4).Value = 3.06
$ws.Cells.Item(66, 5).Value = "202507021326"

$ws.Cells.Item(67, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(67, 2).Value = "518880"
$ws.Cells.Item(67, 3).Value = "黄金ETF"
$ws.Cells.Item(67, 4).Value = 4.88
$ws.Cells.Item(67, 5).Value = "202507021326"

$ws.Cells.Item(68, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(68, 2).Value = "510300"
$ws.Cells.Item(68, 3).Value = "沪深300ETF"
$ws.Cells.Item(68, 4).Value = 5.01
$ws.Cells.Item(68, 5).Value = "202507021326"

$ws.Cells.Item(69, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(69, 2).Value = "HK02899"
$ws.Cells.Item(69, 3).Value = "紫金矿业"
$ws.Cells.Item(69, 4).Value = 9.89
$ws.Cells.Item(69, 5).Value = "202507021326"

$ws.Cells.Item(70, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(70, 2).Value = "HK06881"
$ws.Cells.Item(70, 3).Value = "中国银河"
$ws.Cells.Item(70, 4).Value = 5.22
$ws.Cells.Item(70, 5).Value = "202507021326"

$ws.Cells.Item(71, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(71, 2).Value = "510050"
$ws.Cells.Item(71, 3).Value = "上证50ETF"
$ws.Cells.Item(71, 4).Value = 5.14
$ws.Cells.Item(71, 5).Value = "202507021326"

$ws.Cells.Item(72, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(72, 2).Value = "600085"
$ws.Cells.Item(72, 3).Value = "同仁堂"
$ws.Cells.Item(72, 4).Value = 1.98
$ws.Cells.Item(72, 5).Value = "202507021326"

$ws.Cells.Item(73, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(73, 2).Value = "600989"
$ws.Cells.Item(73, 3).Value = "宝丰能源"
$ws.Cells.Item(73, 4).Value = 1
$ws.Cells.Item(73, 5).Value = "202507021326"

$ws.Cells.Item(74, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(74, 2).Value = "601899"
$ws.Cells.Item(74, 3).Value = "紫金矿业"
$ws.Cells.Item(74, 4).Value = 9.92
$ws.Cells.Item(74, 5).Value = "202507021326"

$ws.Cells.Item(75, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(75, 2).Value = "601688"
$ws.Cells.Item(75, 3).Value = "华泰证券"
$ws.Cells.Item(75, 4).Value = 5
$ws.Cells.Item(75, 5).Value = "202507021326"

$ws.Cells.Item(76, 1).Value = "大智 (稳健智远)"
$ws.Cells.Item(76, 2).Value = "600380"
$ws.Cells.Item(76, 3).Value = "健康元"
$ws.Cells.Item(76, 4).Value = 10
$ws.Cells.Item(76, 5).Value = "202507021326"

# Sheet 2: 大成投资组合 - append daily snapshot rows 35-48
$ws = $wb.Worksheets("大成投资组合")

$ws.Range("B35:B48").NumberFormat = "@"
$ws.Range("E35:E48").NumberFormat = "@"

$ws.Cells.Item(35, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(35, 2).Value = "603119"
$ws.Cells.Item(35, 3).Value = "浙江荣泰"
$ws.Cells.Item(35, 4).Value = 42.95
$ws.Cells.Item(35, 5).Value = "202507021030"

$ws.Cells.Item(36, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(36, 2).Value = "688290"
$ws.Cells.Item(36, 3).Value = "景业智能"
$ws.Cells.Item(36, 4).Value = 7.46
$ws.Cells.Item(36, 5).Value = "202507021030"

$ws.Cells.Item(37, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(37, 2).Value = "513100"
$ws.Cells.Item(37, 3).Value = "纳指ETF"
$ws.Cells.Item(37, 4).Value = 4.85
$ws.Cells.Item(37, 5).Value = "202507021030"

$ws.Cells.Item(38, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(38, 2).Value = "513290"
$ws.Cells.Item(38, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(38, 4).Value = 0.93
$ws.Cells.Item(38, 5).Value = "202507021030"

$ws.Cells.Item(39, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(39, 2).Value = "000725"
$ws.Cells.Item(39, 3).Value = "京东方A"
$ws.Cells.Item(39, 4).Value = 4.84
$ws.Cells.Item(39, 5).Value = "202507021030"

$ws.Cells.Item(40, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(40, 2).Value = "159781"
$ws.Cells.Item(40, 3).Value = "科创创业ETF"
$ws.Cells.Item(40, 4).Value = 9.93
$ws.Cells.Item(40, 5).Value = "202507021030"

$ws.Cells.Item(41, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(41, 2).Value = "HK01896"
$ws.Cells.Item(41, 3).Value = "猫眼娱乐"
$ws.Cells.Item(41, 4).Value = 0.97
$ws.Cells.Item(41, 5).Value = "202507021030"

$ws.Cells.Item(42, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(42, 2).Value = "603119"
$ws.Cells.Item(42, 3).Value = "浙江荣泰"
$ws.Cells.Item(42, 4).Value = 42.95
$ws.Cells.Item(42, 5).Value = "202507021326"

$ws.Cells.Item(43, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(43, 2).Value = "688290"
$ws.Cells.Item(43, 3).Value = "景业智能"
$ws.Cells.Item(43, 4).Value = 7.46
$ws.Cells.Item(43, 5).Value = "202507021326"

$ws.Cells.Item(44, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(44, 2).Value = "513100"
$ws.Cells.Item(44, 3).Value = "纳指ETF"
$ws.Cells.Item(44, 4).Value = 4.85
$ws.Cells.Item(44, 5).Value = "202507021326"

$ws.Cells.Item(45, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(45, 2).Value = "513290"
$ws.Cells.Item(45, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(45, 4).Value = 0.93
$ws.Cells.Item(45, 5).Value = "202507021326"

$ws.Cells.Item(46, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(46, 2).Value = "000725"
$ws.Cells.Item(46, 3).Value = "京东方A"
$ws.Cells.Item(46, 4).Value = 4.84
$ws.Cells.Item(46, 5).Value = "202507021326"

$ws.Cells.Item(47, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(47, 2).Value = "159781"
$ws.Cells.Item(47, 3).Value = "科创创业ETF"
$ws.Cells.Item(47, 4).Value = 9.93
$ws.Cells.Item(47, 5).Value = "202507021326"

$ws.Cells.Item(48, 1).Value = "大成 (锐进先锋)"
$ws.Cells.Item(48, 2).Value = "HK01896"
$ws.Cells.Item(48, 3).Value = "猫眼娱乐"
$ws.Cells.Item(48, 4).Value = 0.97
$ws.Cells.Item(48, 5).Value = "202507021326"

# Sheet 3: 我的投资组合 - append daily snapshot rows 81-115
$ws = $wb.Worksheets("我的投资组合")

$ws.Range("B81:B115").NumberFormat = "@"
$ws.Range("G81:G115").NumberFormat = "@"

$ws.Cells.Item(81, 1).Value = "范式进化投资组合"
$ws.Cells.Item(81, 2).Value = "600900"
$ws.Cells.Item(81, 3).Value = "长江电力"
$ws.Cells.Item(81, 4).Value = "大智"
$ws.Cells.Item(81, 5).Value = 30.04
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = "202507021030"

$ws.Cells.Item(82, 1).Value = "范式进化投资组合"
$ws.Cells.Item(82, 2).Value = "000333"
$ws.Cells.Item(82, 3).Value = "美的集团"
$ws.Cells.Item(82, 4).Value = "大智"
$ws.Cells.Item(82, 5).Value = 3.06
$ws.Cells.Item(82, 6).Value = 1.02
$ws.Cells.Item(82, 7).Value = "202507021030"

$ws.Cells.Item(83, 1).Value = "范式进化投资组合"
$ws.Cells.Item(83, 2).Value = "603119"
$ws.Cells.Item(83, 3).Value = "浙江荣泰"
$ws.Cells.Item(83, 4).Value = "大成"
$ws.Cells.Item(83, 5).Value = 42.8
$ws.Cells.Item(83, 6).Value = 1.14
$ws.Cells.Item(83, 7).Value = "202507021030"

$ws.Cells.Item(84, 1).Value = "范式进化投资组合"
$ws.Cells.Item(84, 2).Value = "518880"
$ws.Cells.Item(84, 3).Value = "黄金ETF"
$ws.Cells.Item(84, 4).Value = "大智"
$ws.Cells.Item(84, 5).Value = 4.87
$ws.Cells.Item(84, 6).Value = 0.98
$ws.Cells.Item(84, 7).Value = "202507021030"

$ws.Cells.Item(85, 1).Value = "范式进化投资组合"
$ws.Cells.Item(85, 2).Value = "510300"
$ws.Cells.Item(85, 3).Value = "沪深300ETF"
$ws.Cells.Item(85, 4).Value = "大智"
$ws.Cells.Item(85, 5).Value = 5.01
$ws.Cells.Item(85, 6).Value = 5
$ws.Cells.Item(85, 7).Value = "202507021030"

$ws.Cells.Item(86, 1).Value = "范式进化投资组合"
$ws.Cells.Item(86, 2).Value = "513100"
$ws.Cells.Item(86, 3).Value = "纳指ETF"
$ws.Cells.Item(86, 4).Value = "大成"
$ws.Cells.Item(86, 5).Value = 4.87
$ws.Cells.Item(86, 6).Value = 1.02
$ws.Cells.Item(86, 7).Value = "202507021030"

$ws.Cells.Item(87, 1).Value = "范式进化投资组合"
$ws.Cells.Item(87, 2).Value = "HK06881"
$ws.Cells.Item(87, 3).Value = "中国银河"
$ws.Cells.Item(87, 4).Value = "大智"
$ws.Cells.Item(87, 5).Value = 5.22
$ws.Cells.Item(87, 6).Value = 1.05
$ws.Cells.Item(87, 7).Value = "202507021030"

$ws.Cells.Item(88, 1).Value = "范式进化投资组合"
$ws.Cells.Item(88, 2).Value = "510050"
$ws.Cells.Item(88, 3).Value = "上证50ETF"
$ws.Cells.Item(88, 4).Value = "大智"
$ws.Cells.Item(88, 5).Value = 5.14
$ws.Cells.Item(88, 6).Value = 5.14
$ws.Cells.Item(88, 7).Value = "202507021030"

$ws.Cells.Item(89, 1).Value = "范式进化投资组合"
$ws.Cells.Item(89, 2).Value = "600085"
$ws.Cells.Item(89, 3).Value = "同仁堂"
$ws.Cells.Item(89, 4).Value = "大智"
$ws.Cells.Item(89, 5).Value = 1.98
$ws.Cells.Item(89, 6).Value = 0.99
$ws.Cells.Item(89, 7).Value = "202507021030"

$ws.Cells.Item(90, 1).Value = "范式进化投资组合"
$ws.Cells.Item(90, 2).Value = "513290"
$ws.Cells.Item(90, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(90, 4).Value = "大成"
$ws.Cells.Item(90, 5).Value = 0.93
$ws.Cells.Item(90, 6).Value = 0.98
$ws.Cells.Item(90, 7).Value = "202507021030"

$ws.Cells.Item(91, 1).Value = "范式进化投资组合"
$ws.Cells.Item(91, 2).Value = "000725"
$ws.Cells.Item(91, 3).Value = "京东方A"
$ws.Cells.Item(91, 4).Value = "大成"
$ws.Cells.Item(91, 5).Value = 4.84
$ws.Cells.Item(91, 6).Value = 5.09
$ws.Cells.Item(91, 7).Value = "202507021030"

$ws.Cells.Item(92, 1).Value = "范式进化投资组合"
$ws.Cells.Item(92, 2).Value = "159781"
$ws.Cells.Item(92, 3).Value = "科创创业ETF"
$ws.Cells.Item(92, 4).Value = "大成"
$ws.Cells.Item(92, 5).Value = 9.96
$ws.Cells.Item(92, 6).Value = 5.25
$ws.Cells.Item(92, 7).Value = "202507021030"

$ws.Cells.Item(93, 1).Value = "范式进化投资组合"
$ws.Cells.Item(93, 2).Value = "600989"
$ws.Cells.Item(93, 3).Value = "宝丰能源"
$ws.Cells.Item(93, 4).Value = "大智"
$ws.Cells.Item(93, 5).Value = 4.8
$ws.Cells.Item(93, 6).Value = 1
$ws.Cells.Item(93, 7).Value = "202507021030"

$ws.Cells.Item(94, 1).Value = "范式进化投资组合"
$ws.Cells.Item(94, 2).Value = "601899"
$ws.Cells.Item(94, 3).Value = "紫金矿业"
$ws.Cells.Item(94, 4).Value = "大智"
$ws.Cells.Item(94, 5).Value = 9.91
$ws.Cells.Item(94, 6).Value = 9.91
$ws.Cells.Item(94, 7).Value = "202507021030"

$ws.Cells.Item(95, 1).Value = "范式进化投资组合"
$ws.Cells.Item(95, 2).Value = "HK02899"
$ws.Cells.Item(95, 3).Value = "紫金矿业"
$ws.Cells.Item(95, 4).Value = "大智"
$ws.Cells.Item(95, 5).Value = 9.89
$ws.Cells.Item(95, 6).Value = 1.11
$ws.Cells.Item(95, 7).Value = "202507021030"

$ws.Cells.Item(96, 1).Value = "范式进化投资组合"
$ws.Cells.Item(96, 2).Value = "688290"
$ws.Cells.Item(96, 3).Value = "景业智能"
$ws.Cells.Item(96, 4).Value = "大成"
$ws.Cells.Item(96, 5).Value = 7.46
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = "202507021030"

$ws.Cells.Item(97, 1).Value = "范式进化投资组合"
$ws.Cells.Item(97, 2).Value = "HK01896"
$ws.Cells.Item(97, 3).Value = "猫眼娱乐"
$ws.Cells.Item(97, 4).Value = "大成"
$ws.Cells.Item(97, 5).Value = 0.97
$ws.Cells.Item(97, 6).Value = 0.2
$ws.Cells.Item(97, 7).Value = "202507021030"

$ws.Cells.Item(98, 1).Value = "范式进化投资组合"
$ws.Cells.Item(98, 2).Value = "600900"
$ws.Cells.Item(98, 3).Value = "长江电力"
$ws.Cells.Item(98, 4).Value = "大智"
$ws.Cells.Item(98, 5).Value = 17
$ws.Cells.Item(98, 6).Value = 1
$ws.Cells.Item(98, 7).Value = "202507021326"

$ws.Cells.Item(99, 1).Value = "范式进化投资组合"
$ws.Cells.Item(99, 2).Value = "000333"
$ws.Cells.Item(99, 3).Value = "美的集团"
$ws.Cells.Item(99, 4).Value = "大智"
$ws.Cells.Item(99, 5).Value = 3.06
$ws.Cells.Item(99, 6).Value = 1.02
$ws.Cells.Item(99, 7).Value = "202507021326"

$ws.Cells.Item(100, 1).Value = "范式进化投资组合"
$ws.Cells.Item(100, 2).Value = "603119"
$ws.Cells.Item(100, 3).Value = "浙江荣泰"
$ws.Cells.Item(100, 4).Value = "大成"
$ws.Cells.Item(100, 5).Value = 42.95
$ws.Cells.Item(100, 6).Value = 1.14
$ws.Cells.Item(100, 7).Value = "202507021326"

$ws.Cells.Item(101, 1).Value = "范式进化投资组合"
$ws.Cells.Item(101, 2).Value = "518880"
$ws.Cells.Item(101, 3).Value = "黄金ETF"
$ws.Cells.Item(101, 4).Value = "大智"
$ws.Cells.Item(101, 5).Value = 4.88
$ws.Cells.Item(101, 6).Value = 0.98
$ws.Cells.Item(101, 7).Value = "202507021326"

$ws.Cells.Item(102, 1).Value = "范式进化投资组合"
$ws.Cells.Item(102, 2).Value = "510300"
$ws.Cells.Item(102, 3).Value = "沪深300ETF"
$ws.Cells.Item(102, 4).Value = "大智"
$ws.Cells.Item(102, 5).Value = 5.01
$ws.Cells.Item(102, 6).Value = 5
$ws.Cells.Item(102, 7).Value = "202507021326"

$ws.Cells.Item(103, 1).Value = "范式进化投资组合"
$ws.Cells.Item(103, 2).Value = "513100"
$ws.Cells.Item(103, 3).Value = "纳指ETF"
$ws.Cells.Item(103, 4).Value = "大成"
$ws.Cells.Item(103, 5).Value = 4.85
$ws.Cells.Item(103, 6).Value = 1.02
$ws.Cells.Item(103, 7).Value = "202507021326"

$ws.Cells.Item(104, 1).Value = "范式进化投资组合"
$ws.Cells.Item(104, 2).Value = "HK06881"
$ws.Cells.Item(104, 3).Value = "中国银河"
$ws.Cells.Item(104, 4).Value = "大智"
$ws.Cells.Item(104, 5).Value = 5.22
$ws.Cells.Item(104, 6).Value = 1.05
$ws.Cells.Item(104, 7).Value = "202507021326"

$ws.Cells.Item(105, 1).Value = "范式进化投资组合"
$ws.Cells.Item(105, 2).Value = "510050"
$ws.Cells.Item(105, 3).Value = "上证50ETF"
$ws.Cells.Item(105, 4).Value = "大智"
$ws.Cells.Item(105, 5).Value = 5.14
$ws.Cells.Item(105, 6).Value = 5.14
$ws.Cells.Item(105, 7).Value = "202507021326"

$ws.Cells.Item(106, 1).Value = "范式进化投资组合"
$ws.Cells.Item(106, 2).Value = "600085"
$ws.Cells.Item(106, 3).Value = "同仁堂"
$ws.Cells.Item(106, 4).Value = "大智"
$ws.Cells.Item(106, 5).Value = 1.98
$ws.Cells.Item(106, 6).Value = 0.99
$ws.Cells.Item(106, 7).Value = "202507021326"

$ws.Cells.Item(107, 1).Value = "范式进化投资组合"
$ws.Cells.Item(107, 2).Value = "513290"
$ws.Cells.Item(107, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(107, 4).Value = "大成"
$ws.Cells.Item(107, 5).Value = 0.93
$ws.Cells.Item(107, 6).Value = 0.98
$ws.Cells.Item(107, 7).Value = "202507021326"

$ws.Cells.Item(108, 1).Value = "范式进化投资组合"
$ws.Cells.Item(108, 2).Value = "000725"
$ws.Cells.Item(108, 3).Value = "京东方A"
$ws.Cells.Item(108, 4).Value = "大成"
$ws.Cells.Item(108, 5).Value = 4.84
$ws.Cells.Item(108, 6).Value = 5.09
$ws.Cells.Item(108, 7).Value = "202507021326"

$ws.Cells.Item(109, 1).Value = "范式进化投资组合"
$ws.Cells.Item(109, 2).Value = "159781"
$ws.Cells.Item(109, 3).Value = "科创创业ETF"
$ws.Cells.Item(109, 4).Value = "大成"
$ws.Cells.Item(109, 5).Value = 9.93
$ws.Cells.Item(109, 6).Value = 5.25
$ws.Cells.Item(109, 7).Value = "202507021326"

$ws.Cells.Item(110, 1).Value = "范式进化投资组合"
$ws.Cells.Item(110, 2).Value = "600989"
$ws.Cells.Item(110, 3).Value = "宝丰能源"
$ws.Cells.Item(110, 4).Value = "大智"
$ws.Cells.Item(110, 5).Value = 1
$ws.Cells.Item(110, 6).Value = 1
$ws.Cells.Item(110, 7).Value = "202507021326"

$ws.Cells.Item(111, 1).Value = "范式进化投资组合"
$ws.Cells.Item(111, 2).Value = "601899"
$ws.Cells.Item(111, 3).Value = "紫金矿业"
$ws.Cells.Item(111, 4).Value = "大智"
$ws.Cells.Item(111, 5).Value = 9.92
$ws.Cells.Item(111, 6).Value = 9.91
$ws.Cells.Item(111, 7).Value = "202507021326"

$ws.Cells.Item(112, 1).Value = "范式进化投资组合"
$ws.Cells.Item(112, 2).Value = "HK02899"
$ws.Cells.Item(112, 3).Value = "紫金矿业"
$ws.Cells.Item(112, 4).Value = "大智"
$ws.Cells.Item(112, 5).Value = 9.89
$ws.Cells.Item(112, 6).Value = 1.11
$ws.Cells.Item(112, 7).Value = "202507021326"

$ws.Cells.Item(113, 1).Value = "范式进化投资组合"
$ws.Cells.Item(113, 2).Value = "HK01896"
$ws.Cells.Item(113, 3).Value = "猫眼娱乐"
$ws.Cells.Item(113, 4).Value = "大成"
$ws.Cells.Item(113, 5).Value = 0.97
$ws.Cells.Item(113, 6).Value = 0.2
$ws.Cells.Item(113, 7).Value = "202507021326"

$ws.Cells.Item(114, 1).Value = "范式进化投资组合"
$ws.Cells.Item(114, 2).Value = "601688"
$ws.Cells.Item(114, 3).Value = "华泰证券"
$ws.Cells.Item(114, 4).Value = "大智"
$ws.Cells.Item(114, 5).Value = 5
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = "202507021326"

$ws.Cells.Item(115, 1).Value = "范式进化投资组合"
$ws.Cells.Item(115, 2).Value = "600380"
$ws.Cells.Item(115, 3).Value = "健康元"
$ws.Cells.Item(115, 4).Value = "大智"
$ws.Cells.Item(115, 5).Value = 10
$ws.Cells.Item(115, 6).Value = 5
$ws.Cells.Item(115, 7).Value = "202507021326"

